# Regenerate all penyata to follow new data and format
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Relabel the revision-tracking rows ("Kali ..." -> "Semakan Kali ...")
#    These labels are repeated for each of the three sections (Merit
#    Pendahuluan, Laporan Atas Talian, JPPM/JDM/JDRM), so every occurrence
#    is updated together to keep the shared text consistent.
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = "Semakan Kali Pertama"
$ws.Range("C22").Value = "Semakan Kali Pertama"
$ws.Range("C28").Value = "Semakan Kali Pertama"

$ws.Range("C17").Value = "Semakan Kali Kedua"
$ws.Range("C23").Value = "Semakan Kali Kedua"
$ws.Range("C29").Value = "Semakan Kali Kedua"

$ws.Range("C18").Value = "Semakan Kali Ketiga"
$ws.Range("C24").Value = "Semakan Kali Ketiga"
$ws.Range("C30").Value = "Semakan Kali Ketiga"

$ws.Range("C19").Value = "Semakan Kali Keempat"
$ws.Range("C25").Value = "Semakan Kali Keempat"
$ws.Range("C31").Value = "Semakan Kali Keempat"

# ---------------------------------------------------------------------------
# 2. Fix the competition-entry names casing (ALL CAPS -> Title Case)
# ---------------------------------------------------------------------------
$ws.Range("C34").Value = "Bouquet Kreatif"
$ws.Range("C35").Value = "Tik Tok Raya"
$ws.Range("C36").Value = "Riang Ria Kuih Raya"
$ws.Range("C37").Value = "Creative Collage"

# ---------------------------------------------------------------------------
# 3. Updated merit/demerit figures for "Penandaan Fail - Kali Ketiga"
#    and demerit figure for "Laporan Atas Talian - Kali Kedua"
# ---------------------------------------------------------------------------
$ws.Range("D18").Value = 3780
$ws.Range("E18").Value = 1970
$ws.Range("E29").Value = 1500

# ---------------------------------------------------------------------------
# 4. Move the "STATEMENT OF HOMEROOM ACCOUNT" title from E4 into D4 and
#    widen its merge so it spans the full table width (D4:G4).
# ---------------------------------------------------------------------------
$title = $ws.Range("E4").Value()
$ws.Range("E4").ClearContents()
$ws.Range("D4").Value = $title
$ws.Range("D4:G4").Merge()

# ---------------------------------------------------------------------------
# 5. Merged-cell layout adjustments to match the refreshed form layout.
# ---------------------------------------------------------------------------
$ws.Range("B5:C5").Merge()
$ws.Range("B12:F12").Merge()
$ws.Range("B15:C15").UnMerge()

$ws.Range("B21:C21").UnMerge()
$ws.Range("B21:E21").Merge()

$ws.Range("B27:C27").UnMerge()
$ws.Range("B27:E27").Merge()

$ws.Range("B33:C33").UnMerge()
$ws.Range("B33:E33").Merge()

$ws.Range("B43:E43").Merge()

# ---------------------------------------------------------------------------
# 6. Reposition/resize the letterhead logo picture.
# ---------------------------------------------------------------------------
$shp = $ws.Shapes.Item(1)
$colLeft = $ws.Cells.Item(1, 2).Left()
$rowTop = $ws.Cells.Item(1, 1).Top()
$shp.Left = $colLeft + 12
$shp.Top = $rowTop + 14.25
$shp.Width = 46.5
$shp.Height = 47.25

# ---------------------------------------------------------------------------
# 7. Print / page setup: fit to one page, centre horizontally, drop the
#    header/footer margins.
# ---------------------------------------------------------------------------
$ps = $ws.PageSetup
$ps.FitToPagesWide = 1
$ps.FitToPagesTall = 1
$ps.CenterHorizontally = $true
$ps.HeaderMargin = 0
$ps.FooterMargin = 0
